$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.41996
$ws.Range("H2").Value = 88.25988000000001
$ws.Range("I2").Value = 0.6207199949605289
$ws.Range("J2").Value = 0.66829493802317
$ws.Range("M2").Value = 1.684496
$ws.Range("N2").Value = 5.053488
$ws.Range("O2").Value = 0.6423607101334534
$ws.Range("P2").Value = 0.7282461611889918
$ws.Range("Q2").Value = 49.55780494016
$ws.Range("R2").Value = 446.02024446144
$ws.Range("S2").Value = 0.398726136756879
$ws.Range("T2").Value = 0.4866832231574087
$ws.Range("G3").Value = 29.41996
$ws.Range("H3").Value = 88.25988000000001
$ws.Range("I3").Value = 0.6207199949605289
$ws.Range("J3").Value = 0.66829493802317
$ws.Range("M3").Value = 0.9277985
$ws.Range("N3").Value = 1.855597
$ws.Range("O3").Value = 0.3538039290807178
$ws.Range("P3").Value = 0.2674056793968462
$ws.Range("Q3").Value = 27.29579475806
$ws.Range("R3").Value = 163.77476854836
$ws.Range("S3").Value = 0.2196131730759985
$ws.Range("T3").Value = 0.178705861939559
$ws.Range("G4").Value = 29.41996
$ws.Range("H4").Value = 88.25988000000001
$ws.Range("I4").Value = 0.6207199949605289
$ws.Range("J4").Value = 0.66829493802317
$ws.Range("M4").Value = 0.01005766666666667
$ws.Range("N4").Value = 0.030173
$ws.Range("O4").Value = 0.003835360785828855
$ws.Range("P4").Value = 0.004348159414162149
$ws.Range("Q4").Value = 0.2958961510266667
$ws.Range("R4").Value = 2.66306535924
$ws.Range("S4").Value = 0.002380685127651497
$ws.Range("T4").Value = 0.002905852926202356
$ws.Range("I5").Value = 0.1515698101047853
$ws.Range("J5").Value = 0.1631868437822795
$ws.Range("M5").Value = 1.684496
$ws.Range("N5").Value = 5.053488
$ws.Range("O5").Value = 0.6423607101334534
$ws.Range("P5").Value = 0.7282461611889918
$ws.Range("Q5").Value = 12.10121656298133
$ws.Range("R5").Value = 108.910949066832
$ws.Range("S5").Value = 0.09736249085370255
$ws.Range("T5").Value = 0.1188401925409927
$ws.Range("I6").Value = 0.1515698101047853
$ws.Range("J6").Value = 0.1631868437822795
$ws.Range("M6").Value = 0.9277985
$ws.Range("N6").Value = 1.855597
$ws.Range("O6").Value = 0.3538039290807178
$ws.Range("P6").Value = 0.2674056793968462
$ws.Range("Q6").Value = 6.665192778913833
$ws.Range("R6").Value = 39.991156673483
$ws.Range("S6").Value = 0.05362599434509131
$ws.Range("T6").Value = 0.04363708883022746
$ws.Range("I7").Value = 0.1515698101047853
$ws.Range("J7").Value = 0.1631868437822795
$ws.Range("M7").Value = 0.01005766666666667
$ws.Range("N7").Value = 0.030173
$ws.Range("O7").Value = 0.003835360785828855
$ws.Range("P7").Value = 0.004348159414162149
$ws.Range("Q7").Value = 0.07225306706077778
$ws.Range("R7").Value = 0.650277603547
$ws.Range("S7").Value = 0.0005813249059914195
$ws.Range("T7").Value = 0.0007095624110593265
$ws.Range("G8").Value = 0.3873096666666667
$ws.Range("H8").Value = 1.161929
$ws.Range("I8").Value = 0.008171692087327698
$ws.Range("J8").Value = 0.008798009571759262
$ws.Range("M8").Value = 1.684496
$ws.Range("N8").Value = 5.053488
$ws.Range("O8").Value = 0.6423607101334534
$ws.Range("P8").Value = 0.7282461611889918
$ws.Range("Q8").Value = 0.6524215842613333
$ws.Range("R8").Value = 5.871794258352
$ws.Range("S8").Value = 0.005249173932207743
$ws.Range("T8").Value = 0.006407116696737688
$ws.Range("G9").Value = 0.3873096666666667
$ws.Range("H9").Value = 1.161929
$ws.Range("I9").Value = 0.008171692087327698
$ws.Range("J9").Value = 0.008798009571759262
$ws.Range("M9").Value = 0.9277985
$ws.Range("N9").Value = 1.855597
$ws.Range("O9").Value = 0.3538039290807178
$ws.Range("P9").Value = 0.2674056793968462
$ws.Range("Q9").Value = 0.3593453277688333
$ws.Range("R9").Value = 2.156071966613
$ws.Range("S9").Value = 0.002891176767734352
$ws.Range("T9").Value = 0.002352637726876242
$ws.Range("G10").Value = 0.3873096666666667
$ws.Range("H10").Value = 1.161929
$ws.Range("I10").Value = 0.008171692087327698
$ws.Range("J10").Value = 0.008798009571759262
$ws.Range("M10").Value = 0.01005766666666667
$ws.Range("N10").Value = 0.030173
$ws.Range("O10").Value = 0.003835360785828855
$ws.Range("P10").Value = 0.004348159414162149
$ws.Range("Q10").Value = 0.003895431524111111
$ws.Range("R10").Value = 0.035058883717
$ws.Range("S10").Value = 0.000031341387385604602567462651
$ws.Range("T10").Value = 0.000038255148145333727090471737
$ws.Range("G11").Value = 10.122265
$ws.Range("H11").Value = 20.24453
$ws.Range("I11").Value = 0.2135656295858028
$ws.Range("J11").Value = 0.153289545846405
$ws.Range("M11").Value = 1.684496
$ws.Range("N11").Value = 5.053488
$ws.Range("O11").Value = 0.6423607101334534
$ws.Range("P11").Value = 0.7282461611889918
$ws.Range("Q11").Value = 17.05091490344
$ws.Range("R11").Value = 102.30548942064
$ws.Range("S11").Value = 0.1371861694808344
$ws.Range("T11").Value = 0.1116325233130484
$ws.Range("G12").Value = 10.122265
$ws.Range("H12").Value = 20.24453
$ws.Range("I12").Value = 0.2135656295858028
$ws.Range("J12").Value = 0.153289545846405
$ws.Range("M12").Value = 0.9277985
$ws.Range("N12").Value = 1.855597
$ws.Range("O12").Value = 0.3538039290807178
$ws.Range("P12").Value = 0.2674056793968462
$ws.Range("Q12").Value = 9.391422283602498
$ws.Range("R12").Value = 37.56568913440999
$ws.Range("S12").Value = 0.07556035886405424
$ws.Range("T12").Value = 0.04099049515149194
$ws.Range("G13").Value = 10.122265
$ws.Range("H13").Value = 20.24453
$ws.Range("I13").Value = 0.2135656295858028
$ws.Range("J13").Value = 0.153289545846405
$ws.Range("M13").Value = 0.01005766666666667
$ws.Range("N13").Value = 0.030173
$ws.Range("O13").Value = 0.003835360785828855
$ws.Range("P13").Value = 0.004348159414162149
$ws.Range("Q13").Value = 0.1018063672816666
$ws.Range("R13").Value = 0.6108382036899999
$ws.Range("S13").Value = 0.0008191012409142389
$ws.Range("T13").Value = 0.0006665273818646862
$ws.Range("G14").Value = 0.2830933333333334
$ws.Range("H14").Value = 0.84928
$ws.Range("I14").Value = 0.005972873261555284
$ws.Range("J14").Value = 0.006430662776386256
$ws.Range("M14").Value = 1.684496
$ws.Range("N14").Value = 5.053488
$ws.Range("O14").Value = 0.6423607101334534
$ws.Range("P14").Value = 0.7282461611889918
$ws.Range("Q14").Value = 0.4768695876266667
$ws.Range("R14").Value = 4.29182628864
$ws.Range("S14").Value = 0.003836739109829768
$ws.Range("T14").Value = 0.004683105480804235
$ws.Range("G15").Value = 0.2830933333333334
$ws.Range("H15").Value = 0.84928
$ws.Range("I15").Value = 0.005972873261555284
$ws.Range("J15").Value = 0.006430662776386256
$ws.Range("M15").Value = 0.9277985
$ws.Range("N15").Value = 1.855597
$ws.Range("O15").Value = 0.3538039290807178
$ws.Range("P15").Value = 0.2674056793968462
$ws.Range("Q15").Value = 0.2626535700266667
$ws.Range("R15").Value = 1.57592142016
$ws.Range("S15").Value = 0.002113226027839421
$ws.Range("T15").Value = 0.001719595748691576
$ws.Range("G16").Value = 0.2830933333333334
$ws.Range("H16").Value = 0.84928
$ws.Range("I16").Value = 0.005972873261555284
$ws.Range("J16").Value = 0.006430662776386256
$ws.Range("M16").Value = 0.01005766666666667
$ws.Range("N16").Value = 0.030173
$ws.Range("O16").Value = 0.003835360785828855
$ws.Range("P16").Value = 0.004348159414162149
$ws.Range("Q16").Value = 0.002847258382222222
$ws.Range("R16").Value = 0.02562532544
$ws.Range("S16").Value = 0.000022908123886094828441928328
$ws.Range("T16").Value = 0.000027961546890446000454902298
